# "new pagination from DB directly"
# Adds a new "issue" tracking column (I) to the feature sheet, records two
# new sub-feature rows (pagination moved from collections.sort/JDBC order by
# to DB-side paging) and marks the earlier "open" items as "done", cross
# linking them to the new issue log entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New strings get interned into the shared-string table in this exact
# order so they land at the same indices the original edit produced
# (55: collections.sort..., 56: issue, 57: cannot do it , error, 58: suspend)
$ws.Range("G6").Value = "collections.sort or JDBC order by"
$ws.Range("I1").Value = "issue"
$ws.Range("I37").Value = "cannot do it , error"
$ws.Range("E17").Value = "suspend"

# --- reuse of already-existing strings (de-duped automatically) ---
$ws.Range("I17").Value = "issue ,can not populate order data into db"
$ws.Range("G26").Value = "collections.sort or JDBC order by"
$ws.Range("E37").Value = "suspend"

# --- Row 6: pagination sub-feature flips from "open" to "done" and now has
# a completion date + an issue note; pick up the same look (fill/number
# format) already used a few rows down for the analogous "done" rows.
$ws.Range("E7").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("E6").Value = "done"
$ws.Range("F6").Value = 43986

# --- Row 26: same pagination sub-feature recorded lower in the sheet ---
$ws.Range("E26").Value = "done"
$ws.Range("F26").Value = 43986

# --- Row 17 & 37: style the new issue / suspend notes like their neighbors
$ws.Range("E18").Copy()
$ws.Range("I17").PasteSpecial(-4122)
$ws.Range("I17").Value = "issue ,can not populate order data into db"

$ws.Range("E18").Copy()
$ws.Range("E37").PasteSpecial(-4122)
$ws.Range("E37").Value = "suspend"

$ws.Range("E9").Select()
